# Edit summary (per the target XML diff):
#   1. Three tables (on slides 14, 15 and 16) switch from the custom
#      table style {9DDE1D7C-C070-499C-BC8C-F6B0B0C3399D} to the
#      built-in table style {AC86114E-3611-4BBF-829E-A40F3C6F88AC}.
#   2. The slide master's theme color scheme (the "Integral" design's
#      "Red Violet" palette, stored in ppt/theme/theme1.xml) is swapped
#      for the standard "Office" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables that used the custom table style.
# ---------------------------------------------------------------------
$newTableStyle = "{AC86114E-3611-4BBF-829E-A40F3C6F88AC}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the slide master's theme colours from "Red Violet" to the
#    standard Office palette (dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink - in that order).
# ---------------------------------------------------------------------
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $hex = $officeColors[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbLong = $b * 65536 + $g * 256 + $r
    $colorScheme.Item($i + 1).RGB = $rgbLong
}
